# Config.xlsx update:
# "added creating report at the beginning of the process and setting asset
#  to be used in producer process"
#
# - Settings sheet gains 4 new rows describing the report template path,
#   the report filename pattern, the asset name used to share the current
#   report filepath between dispatcher/producers, and the orchestrator
#   asset folder.
# - Settings becomes the active/selected sheet (instead of Constants).
# - Selection state on each sheet moves to reflect where the author ended
#   up after editing.
# - A trailing, data-less formatting-only row at the bottom of the Assets
#   sheet is dropped.

$wb = $excel.ActiveWorkbook
$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets    = $wb.Worksheets.Item("Assets")

# --- Settings sheet: new rows 9-12 -----------------------------------
# (Values are written in this particular order so that new shared-string
# entries are registered in the same sequence the original authored file
# used - it does not affect the visible result, only string-table order.)

$settings.Range("A9").Value = "ReportTemplateFilepath"
$settings.Range("B9").Value = "Data\ReportTemplate.xlsx"
$settings.Range("C9").Value = "path to file which contains template for output report"

$settings.Range("A10").Value = "ReportFilepathPattern"

$settings.Range("B11").Value = "CurrentReportFilepath"
$settings.Range("B10").Value = "D:\UiPath Projects\ListUSAPrivateSchools\Reports\USAPrivateSchools_<replace_timestamp>.xlsx"
$settings.Range("A11").Value = "CurrentReportFilepathAssetName"

$settings.Range("C10").Value = "path to file which contains template for output report"
$settings.Range("C11").Value = "Asset to store filepath to report created by dispatcher - will be utilised by one or many producers to populate the single report"

$settings.Range("A12").Value = "OrchestratorAssetFolder"
$settings.Range("B12").Value = "ListUSAPrivateSchools"
$settings.Range("C12").Value = "Name of the orchestrator folder which contains Asset to store current report filepath"

# --- Assets sheet: drop the unused formatting-only row at the bottom --
$assets.Rows.Item(1000).Delete()

# --- Selection / active-sheet bookkeeping -----------------------------
$assets.Range("A13").Select() | Out-Null
$constants.Range("A19").Select() | Out-Null

$settings.Activate()
$settings.Range("C31").Select() | Out-Null
